$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.429.97"
$ws.Range("E2").Value = "  +0.96%  "
$ws.Range("D3").Value = "2.592.48"
$ws.Range("E3").Value = "  -0.21%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'572.35"
$ws.Range("E5").Value = "  +3.19%  "
$ws.Range("D6").Value = "'144.32"
$ws.Range("E6").Value = "  +0.91%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "'0.599"
$ws.Range("E8").Value = "  -0.19%  "
$ws.Range("D9").Value = "2.605.22"
$ws.Range("E9").Value = "  -0.13%  "
$ws.Range("D10").Value = "'6.69"
$ws.Range("E10").Value = "  -1.29%  "
$ws.Range("E11").Value = "  +3.81%  "
$ws.Range("D12").Value = "'0.157"
$ws.Range("E12").Value = "  +10.61%  "
$ws.Range("D13").Value = "'0.348"
$ws.Range("E13").Value = "  +3.64%  "
$ws.Range("D14").Value = "3.047.69"
$ws.Range("E14").Value = "  -0.24%  "
$ws.Range("D15").Value = "59.387.96"
$ws.Range("E15").Value = "  +0.94%  "
$ws.Range("D16").Value = "'22.61"
$ws.Range("E16").Value = "  +8.50%  "
$ws.Range("E17").Value = "  +4.23%  "
$ws.Range("D18").Value = "2.596.80"
$ws.Range("D19").Value = "'4.56"
$ws.Range("E19").Value = "  +1.84%  "
$ws.Range("D20").Value = "'337.45"
$ws.Range("E20").Value = "  -0.18%  "
$ws.Range("E21").Value = "  +1.71%  "
$ws.Range("D22").Value = "'6.23"
$ws.Range("E22").Value = "  +0.80%  "
$ws.Range("D24").Value = "'64.50"
$ws.Range("E24").Value = "  -3.36%  "
$ws.Range("D25").Value = "'0.458"
$ws.Range("E25").Value = "  +6.76%  "
$ws.Range("B26").Value = "Binance-PegBSC-USD"
$ws.Range("C26").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  +0.28%  "
$ws.Range("B27").Value = "Kaspa"
$ws.Range("C27").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D27").Value = "'0.161"
$ws.Range("E27").Value = "  -0.09%  "
$ws.Range("E28").Value = "  +1.98%  "
$ws.Range("D29").Value = "0.0₃0785"
$ws.Range("E29").Value = "  +4.00%  "
$ws.Range("E30").Value = "  +0.01%  "
$ws.Range("E31").Value = "  +0.75%  "
$ws.Range("D32").Value = "'6.11"
$ws.Range("E32").Value = "  +1.10%  "
$ws.Range("D33").Value = "'159.06"
$ws.Range("E33").Value = "  +3.04%  "
$ws.Range("D34").Value = "'19.08"
$ws.Range("E34").Value = "  +0.51%  "
$ws.Range("D35").Value = "'4.06"
$ws.Range("E35").Value = "  +3.15%  "
$ws.Range("E36").Value = "  +1.79%  "
$ws.Range("D37").Value = "'0.888"
$ws.Range("E37").Value = "  +3.61%  "
$ws.Range("D38").Value = "'0.882"
$ws.Range("E38").Value = "  -1.86%  "
$ws.Range("E39").Value = "  +2.05%  "
$ws.Range("D40").Value = "'37.17"
$ws.Range("E40").Value = "  +0.91%  "
$ws.Range("D41").Value = "'296.37"
$ws.Range("E41").Value = "  +4.43%  "
$ws.Range("E42").Value = "  +2.23%  "
$ws.Range("E43").Value = "  +0.02%  "
$ws.Range("D44").Value = "'0.0979"
$ws.Range("E44").Value = "  +2.35%  "
$ws.Range("E45").Value = "  -0.45%  "
$ws.Range("D46").Value = "'0.0539"
$ws.Range("E46").Value = "  +0.86%  "
$ws.Range("D47").Value = "'19.31"
$ws.Range("E47").Value = "  +2.83%  "
$ws.Range("D48").Value = "'10.65"
$ws.Range("E48").Value = "  +0.13%  "
$ws.Range("E49").Value = "  +2.53%  "
$ws.Range("D50").Value = "'124.56"
$ws.Range("E50").Value = "  +5.53%  "
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").Value = "'18.67"
$ws.Range("E51").Value = "  +3.74%  "
